# Regenerate orders with updated distance/size codes.
#
# The experiment's distance conditions were renumbered:
#   D64 -> D69, D51 -> D55, D80 -> D86
# and the "large" size code was renumbered:
#   S30 -> S31
# (S20 and S25 are unchanged.)
#
# These tokens appear, delimited by "_" or "." (or string start/end), inside
# every cell that encodes a trial's condition name, left/right stimulus
# filename, or the standalone Distance/Size lookup columns. We walk the
# sheet's used range and rewrite every text cell whose tokens match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

function Transform($s) {
    $out = $s
    $out = $out -replace '(?<![A-Za-z0-9])D64(?![A-Za-z0-9])', 'D69'
    $out = $out -replace '(?<![A-Za-z0-9])D51(?![A-Za-z0-9])', 'D55'
    $out = $out -replace '(?<![A-Za-z0-9])D80(?![A-Za-z0-9])', 'D86'
    $out = $out -replace '(?<![A-Za-z0-9])S30(?![A-Za-z0-9])', 'S31'
    return $out
}

$changed = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $txt = $cell.Text
        if ($txt -eq $null) {
            continue
        }
        $new = Transform $txt
        if ($new -ne $txt) {
            $cell.Value = $new
            $changed = $changed + 1
        }
    }
}
